# Adjusted risk calc formula
# Update the "Total Risk" and "Current Risk" columns for the four
# watershed risk rows that currently read Rank=3 / Total Risk=4 / Current Risk=L
# (LF1, LF9, LF11, LF39) so that Total Risk becomes 6 and Current Risk becomes M.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $rank = $t.Cell($r, 3).Range.Text.Trim([char]13, [char]7)
    $total = $t.Cell($r, 4).Range.Text.Trim([char]13, [char]7)
    $current = $t.Cell($r, 5).Range.Text.Trim([char]13, [char]7)

    if ($rank -eq "3" -and $total -eq "4" -and $current -eq "L") {
        $totalRange = $t.Cell($r, 4).Range
        [void]$totalRange.MoveEnd(1, -1)
        $totalRange.Text = "6"

        $currentRange = $t.Cell($r, 5).Range
        [void]$currentRange.MoveEnd(1, -1)
        $currentRange.Text = "M"
    }
}
